$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 210.25
$ws.Range("I5").Value = 210.25
$ws.Range("K5").Value = 210.25
$ws.Range("M5").Value = -95.25
$ws.Range("H33").Value = 297.34616
$ws.Range("I33").Value = 304.84
$ws.Range("K33").Value = 304.84
$ws.Range("M33").Value = -75.83999999999997
$ws.Range("H137").Value = 5081.533
$ws.Range("I137").Value = 1871.2778
$ws.Range("J137").Value = 9896.916999999999
$ws.Range("K137").Value = 5613.8334
$ws.Range("L137").Value = 29690.751
$ws.Range("M137").Value = -3063.8334
$ws.Range("N137").Value = -34790.751
$ws.Range("H138").Value = 4455.4756
$ws.Range("I138").Value = 4847.2856
$ws.Range("J138").Value = 4338.766
$ws.Range("K138").Value = 14541.8568
$ws.Range("L138").Value = 13016.298
$ws.Range("M138").Value = -9401.856800000001
$ws.Range("N138").Value = -23296.298

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4367.8696
$ws.Range("I32").Value = 4430.2954
$ws.Range("K32").Value = 4430.2954
$ws.Range("M32").Value = -4143.2954
$ws.Range("H44").Value = 59999
$ws.Range("J44").Value = 59999
$ws.Range("L44").Value = 59999
$ws.Range("N44").Value = -60975
$ws.Range("H45").Value = 45533764
$ws.Range("I45").Value = 122986.71
$ws.Range("K45").Value = 122986.71
$ws.Range("M45").Value = -122609.71
$ws.Range("H46").Value = 29576
$ws.Range("J46").Value = 29152
$ws.Range("L46").Value = 29152
$ws.Range("N46").Value = -29790
$ws.Range("H55").Value = 50483.668
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 65725.5
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 65725.5
$ws.Range("M55").Value = -19685
$ws.Range("N55").Value = -66355.5
$ws.Range("H102").Value = 1967.875
$ws.Range("I102").Value = 1751.8857
$ws.Range("K102").Value = 1751.8857
$ws.Range("M102").Value = -129.8857
$ws.Range("H110").Value = 1787.2273
$ws.Range("I110").Value = 1758.9474
$ws.Range("J110").Value = 1966.3334
$ws.Range("K110").Value = 1758.9474
$ws.Range("L110").Value = 1966.3334
$ws.Range("M110").Value = 286.0526
$ws.Range("N110").Value = -6056.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("H22").Value = 462.25
$ws.Range("I22").Value = 466.33334
$ws.Range("K22").Value = 466.33334
$ws.Range("M22").Value = -293.33334
$ws.Range("H24").Value = 3700
$ws.Range("I24").Value = 3700
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3700
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3465
$ws.Range("N24").Value = $null
$ws.Range("H99").Value = 5065.3335
$ws.Range("J99").Value = 6998
$ws.Range("L99").Value = 6998
$ws.Range("N99").Value = -9994
$ws.Range("H105").Value = 10401612
$ws.Range("I105").Value = 556545.5600000001
$ws.Range("K105").Value = 556545.5600000001
$ws.Range("M105").Value = -554798.5600000001
$ws.Range("H107").Value = 2960105.8
$ws.Range("I107").Value = 3206551.5
$ws.Range("K107").Value = 3206551.5
$ws.Range("M107").Value = -3204631.5
$ws.Range("H140").Value = 59999
$ws.Range("J140").Value = 59999
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1817.6
$ws.Range("I22").Value = 1970.7142
$ws.Range("K22").Value = 1970.7142
$ws.Range("M22").Value = -1620.7142
$ws.Range("H31").Value = 2162.689
$ws.Range("I31").Value = 1709.2222
$ws.Range("J31").Value = 2465
$ws.Range("K31").Value = 1709.2222
$ws.Range("L31").Value = 2465
$ws.Range("M31").Value = -1414.2222
$ws.Range("N31").Value = -3055
$ws.Range("H34").Value = 2162.689
$ws.Range("I34").Value = 1709.2222
$ws.Range("J34").Value = 2465
$ws.Range("K34").Value = 1709.2222
$ws.Range("L34").Value = 2465
$ws.Range("M34").Value = -1507.2222
$ws.Range("N34").Value = -2869
$ws.Range("H99").Value = 111115450
$ws.Range("I99").Value = 333334660
$ws.Range("J99").Value = 5835.6665
$ws.Range("K99").Value = 333334660
$ws.Range("L99").Value = 5835.6665
$ws.Range("M99").Value = -333333162
$ws.Range("N99").Value = -8831.666499999999
$ws.Range("H126").Value = 111115450
$ws.Range("I126").Value = 333334660
$ws.Range("J126").Value = 5835.6665
$ws.Range("K126").Value = 1000003980
$ws.Range("L126").Value = 17506.9995
$ws.Range("M126").Value = -1000001510
$ws.Range("N126").Value = -22446.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 780.2917
$ws.Range("I5").Value = 558.1818
$ws.Range("J5").Value = 968.2308
$ws.Range("K5").Value = 1674.5454
$ws.Range("L5").Value = 2904.6924
$ws.Range("M5").Value = -1562.5454
$ws.Range("N5").Value = -3128.6924
$ws.Range("H68").Value = 1565759.1
$ws.Range("J68").Value = 1789082.5
$ws.Range("L68").Value = 5367247.5
$ws.Range("N68").Value = -5368869.5
$ws.Range("H71").Value = 1565759.1
$ws.Range("J71").Value = 1789082.5
$ws.Range("L71").Value = 16101742.5
$ws.Range("N71").Value = -16109854.5
$ws.Range("H121").Value = 6747400
$ws.Range("J121").Value = 100828
$ws.Range("L121").Value = 302484
$ws.Range("N121").Value = -305104
$ws.Range("H122").Value = 3036.2727
$ws.Range("J122").Value = 3036.2727
$ws.Range("L122").Value = 27326.4543
$ws.Range("N122").Value = -32226.4543
$ws.Range("H131").Value = 9300.764999999999
$ws.Range("J131").Value = 1525.1538
$ws.Range("L131").Value = 4575.4614
$ws.Range("N131").Value = -14655.4614
$ws.Range("H132").Value = 5775.0293
$ws.Range("I132").Value = 3777.5625
$ws.Range("K132").Value = 33998.0625
$ws.Range("M132").Value = -31468.0625
$ws.Range("H135").Value = 780.2917
$ws.Range("I135").Value = 558.1818
$ws.Range("J135").Value = 968.2308
$ws.Range("K135").Value = 5023.6362
$ws.Range("L135").Value = 8714.0772
$ws.Range("M135").Value = -2488.6362
$ws.Range("N135").Value = -13784.0772

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 10008
$ws.Range("J29").Value = 10008
$ws.Range("L29").Value = 10008
$ws.Range("N29").Value = -10588
$ws.Range("H70").Value = 45469520
$ws.Range("I70").Value = 83338584
$ws.Range("J70").Value = 26649
$ws.Range("K70").Value = 83338584
$ws.Range("L70").Value = 26649
$ws.Range("M70").Value = -83338314
$ws.Range("N70").Value = -27189
$ws.Range("H73").Value = 45469520
$ws.Range("I73").Value = 83338584
$ws.Range("J73").Value = 26649
$ws.Range("K73").Value = 83338584
$ws.Range("L73").Value = 26649
$ws.Range("M73").Value = -83337648
$ws.Range("N73").Value = -28521
$ws.Range("H102").Value = 35722196
$ws.Range("I102").Value = 41670748
$ws.Range("J102").Value = 30875
$ws.Range("K102").Value = 41670748
$ws.Range("L102").Value = 30875
$ws.Range("M102").Value = -41669126
$ws.Range("N102").Value = -34119
$ws.Range("H132").Value = 2120.7334
$ws.Range("I132").Value = 2222.889
$ws.Range("J132").Value = 1967.5
$ws.Range("K132").Value = 6668.667
$ws.Range("L132").Value = 5902.5
$ws.Range("M132").Value = -4138.667
$ws.Range("N132").Value = -10962.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 48679.48
$ws.Range("I40").Value = 77643.625
$ws.Range("K40").Value = 77643.625
$ws.Range("M40").Value = -77507.625
$ws.Range("H46").Value = 5175
$ws.Range("J46").Value = 5266.6665
$ws.Range("L46").Value = 5266.6665
$ws.Range("N46").Value = -5642.6665
$ws.Range("H58").Value = 6826.6665
$ws.Range("J58").Value = 11500
$ws.Range("L58").Value = 11500
$ws.Range("N58").Value = -12020
$ws.Range("H61").Value = 7579.375
$ws.Range("I61").Value = 6866.923
$ws.Range("J61").Value = 10666.667
$ws.Range("K61").Value = 6866.923
$ws.Range("L61").Value = 10666.667
$ws.Range("M61").Value = -6664.923
$ws.Range("N61").Value = -11070.667
$ws.Range("H82").Value = 2300.4
$ws.Range("I82").Value = 1750.5
$ws.Range("K82").Value = 1750.5
$ws.Range("M82").Value = -1389.5
$ws.Range("H85").Value = 2300.4
$ws.Range("I85").Value = 1750.5
$ws.Range("K85").Value = 1750.5
$ws.Range("M85").Value = -502.5
$ws.Range("H93").Value = 1635.6666
$ws.Range("I93").Value = 1561.3846
$ws.Range("J93").Value = 1911.5714
$ws.Range("K93").Value = 1561.3846
$ws.Range("L93").Value = 1911.5714
$ws.Range("M93").Value = -313.3846000000001
$ws.Range("N93").Value = -4407.5714
$ws.Range("H113").Value = 7579.375
$ws.Range("I113").Value = 6866.923
$ws.Range("J113").Value = 10666.667
$ws.Range("K113").Value = 6866.923
$ws.Range("L113").Value = 10666.667
$ws.Range("M113").Value = -4696.923
$ws.Range("N113").Value = -15006.667
$ws.Range("H122").Value = 7777
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -28231
$ws.Range("H134").Value = 109997.5
$ws.Range("J134").Value = 109997.5
$ws.Range("L134").Value = 109997.5
$ws.Range("N134").Value = -120137.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1717.4
$ws.Range("I96").Value = 1200
$ws.Range("J96").Value = 2062.3333
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 2062.3333
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -4808.3333
$ws.Range("H100").Value = 83335540
$ws.Range("J100").Value = 142859500
$ws.Range("L100").Value = 285719000
$ws.Range("N100").Value = -285720082
$ws.Range("H122").Value = 6947500.5
$ws.Range("I122").Value = 3140.4
$ws.Range("J122").Value = 41669300
$ws.Range("K122").Value = 9421.200000000001
$ws.Range("L122").Value = 125007900
$ws.Range("M122").Value = -6971.200000000001
$ws.Range("N122").Value = -125012800
